$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "三花智控"
$ws.Range("B2").Value = "闻泰科技"
$ws.Range("C2").Value = "三花智控"
$ws.Range("B3").Value = "三花智控"
$ws.Range("A4").Value = "中国核建"
$ws.Range("B4").Value = "山子高科"
$ws.Range("C4").Value = "山子高科"
$ws.Range("A5").Value = "闻泰科技"
$ws.Range("B5").Value = "北方稀土"
$ws.Range("C5").Value = "澄星股份"
$ws.Range("A6").Value = "北方稀土"
$ws.Range("B6").Value = "中国核建"
$ws.Range("C6").Value = "中国核建"
$ws.Range("A7").Value = "卧龙电驱"
$ws.Range("B7").Value = "至纯科技"
$ws.Range("C7").Value = "安泰科技"
$ws.Range("A8").Value = "至纯科技"
$ws.Range("B8").Value = "卧龙电驱"
$ws.Range("C8").Value = "上海电力"
$ws.Range("A9").Value = "安泰科技"
$ws.Range("B9").Value = "包钢股份"
$ws.Range("C9").Value = "卧龙电驱"
$ws.Range("A10").Value = "国光连锁"
$ws.Range("B10").Value = "盛和资源"
$ws.Range("C10").Value = "海光信息"
$ws.Range("A11").Value = "合锻智能"
$ws.Range("C11").Value = "至纯科技"
$ws.Range("A12").Value = "合肥城建"
$ws.Range("B12").Value = "工业富联"
$ws.Range("C12").Value = "合肥城建"
$ws.Range("A13").Value = "凯美特气"
$ws.Range("B13").Value = "黄河旋风"
$ws.Range("C13").Value = "北方稀土"
$ws.Range("A14").Value = "上海电力"
$ws.Range("B14").Value = "凯美特气"
$ws.Range("C14").Value = "国电南自"
$ws.Range("A15").Value = "工业富联"
$ws.Range("B15").Value = "五洲新春"
$ws.Range("C15").Value = "凯美特气"
$ws.Range("A16").Value = "阳光电源"
$ws.Range("B16").Value = "长城军工"
$ws.Range("C16").Value = "国光连锁"
$ws.Range("A17").Value = "盛和资源"
$ws.Range("B17").Value = "合肥城建"
$ws.Range("C17").Value = "合锻智能"
$ws.Range("A18").Value = "黄河旋风"
$ws.Range("B18").Value = "紫金矿业"
$ws.Range("C18").Value = "阳光电源"
$ws.Range("A19").Value = "澄星股份"
$ws.Range("B19").Value = "合锻智能"
$ws.Range("C19").Value = "华建集团"
$ws.Range("A20").Value = "长城军工"
$ws.Range("B20").Value = "上海电力"
$ws.Range("C20").Value = "紫金矿业"
$ws.Range("A21").Value = "紫金矿业"
$ws.Range("B21").Value = "国光连锁"
$ws.Range("C21").Value = "工业富联"
